# This script updates the CV results workbook to reflect a re-run of the
# cross-validation analysis (features edit, features crosscor).
#
# Sheets (in workbook order):
#   1 CV_Summary          - per-frequency summary stats (A:H), rows 2-5
#   2 CV_Scores_Detail     - per-fold accuracy scores (A:C), rows 2-21
#   3 Label_Distribution   - per-class sample counts (A:D), rows 2-17
#   4 Analysis_Info        - run metadata (A:B), rows 2-6

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: CV_Summary
# Row order (by Frequency) changes from 30/10/20/40 to 30/40/10/20, and all
# the metric columns (Mean_Accuracy, Std_Accuracy, Best_Fold, Worst_Fold,
# Stability_Score, N_Samples, N_Features) get new values from the re-run.
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("CV_Summary")

$summaryRows = @(
    @{ Row = 2; Frequency = "30hz"; Mean = 0.921875;  Std = 0.0369754986443726; Best = 0.96875;  Worst = 0.859375; Stability = 24.93204377163592; N_Samples = 320; N_Features = 110 },
    @{ Row = 3; Frequency = "40hz"; Mean = 0.88125;   Std = 0.04485218779502289; Best = 0.9375;   Worst = 0.8125;   Stability = 19.64786662244207; N_Samples = 320; N_Features = 110 },
    @{ Row = 4; Frequency = "10hz"; Mean = 0.86875;   Std = 0.02538762001448738; Best = 0.90625;  Worst = 0.84375;  Stability = 34.21942101347193; N_Samples = 320; N_Features = 110 },
    @{ Row = 5; Frequency = "20hz"; Mean = 0.853125;  Std = 0.03365728004459066; Best = 0.890625; Worst = 0.8125;   Stability = 25.34740613013533; N_Samples = 320; N_Features = 110 }
)

foreach ($r in $summaryRows) {
    $wsSummary.Cells.Item($r.Row, 1).Value = $r.Frequency
    $wsSummary.Cells.Item($r.Row, 2).Value = $r.Mean
    $wsSummary.Cells.Item($r.Row, 3).Value = $r.Std
    $wsSummary.Cells.Item($r.Row, 4).Value = $r.Best
    $wsSummary.Cells.Item($r.Row, 5).Value = $r.Worst
    $wsSummary.Cells.Item($r.Row, 6).Value = $r.Stability
    $wsSummary.Cells.Item($r.Row, 7).Value = $r.N_Samples
    $wsSummary.Cells.Item($r.Row, 8).Value = $r.N_Features
}

# ---------------------------------------------------------------------------
# Sheet 2: CV_Scores_Detail
# Row (Frequency, Fold) layout is unchanged; only the Accuracy column (C)
# gets new values from the re-run.
# ---------------------------------------------------------------------------
$wsDetail = $wb.Worksheets.Item("CV_Scores_Detail")

$detailAccuracy = @{
    2  = 0.859375
    3  = 0.90625
    4  = 0.84375
    5  = 0.890625
    6  = 0.84375
    7  = 0.875
    8  = 0.8125
    9  = 0.875
    10 = 0.890625
    11 = 0.8125
    12 = 0.96875
    13 = 0.9375
    14 = 0.9375
    15 = 0.859375
    16 = 0.90625
    17 = 0.875
    18 = 0.921875
    19 = 0.859375
    20 = 0.8125
    21 = 0.9375
}

foreach ($row in $detailAccuracy.Keys) {
    $wsDetail.Cells.Item($row, 3).Value = $detailAccuracy[$row]
}

# ---------------------------------------------------------------------------
# Sheet 3: Label_Distribution
# Every row's Count column (C) doubles from 40 to 80 (Percentage stays 25).
# ---------------------------------------------------------------------------
$wsLabel = $wb.Worksheets.Item("Label_Distribution")

for ($row = 2; $row -le 17; $row++) {
    $wsLabel.Cells.Item($row, 3).Value = 80
}

# ---------------------------------------------------------------------------
# Sheet 4: Analysis_Info
# Worst Performing Frequency changes (40hz -> 20hz) to match the new
# CV_Summary numbers, and the Analysis Date is updated to the new run time.
# ---------------------------------------------------------------------------
$wsInfo = $wb.Worksheets.Item("Analysis_Info")

$wsInfo.Cells.Item(4, 2).Value = "20hz"
$wsInfo.Cells.Item(6, 2).Value = "2025-10-06 15:13:44"
